$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 50 -> 40, Wrong total -12 -> -24, Max text "50 / 140" -> "16 / 112"
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -24
$ws.Range("E12").Value = "16 / 112"
